$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "lat"
$ws.Range("E1").Value = "lon"

$ws.Range("E2").Select()
